# Automatische test-sync: 2025-08-08 20:01:50
# Appends the newly received mail-log entry (row 4) to the "Logs" sheet,
# extends the conditional-formatting ranges to cover it, and refreshes the
# "Dashboard" category tally / ordering to match.

$wb = $excel.ActiveWorkbook

# ---- Logs sheet: append the new mail entry as row 4 -----------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A4").Value = "Kun jij dit even regelen?"
$logs.Range("B4").Value = '"Testbedrijf 123 B.V." <admin@testbedrijf123.nl>'
$logs.Range("C4").Value = "Testmail #1: Kun jij dit even regelen?`nTestbedrijf 123 B.V."
$logs.Range("D4").Value = "Overig"
$logs.Range("E4").Value = "Bedankt, we hebben dit doorgestuurd naar support@testbedrijf123.nl."
$logs.Range("F4").Value = "2025-08-08 20:01:49"
$logs.Range("G4").Value = "Ja"
$logs.Range("H4").Value = "Ja"
$logs.Range("I4").Value = "Nee"
$logs.Range("J4").Value = "Nee"

# Undo the implicit row auto-fit triggered by C4's embedded line break so the
# new row keeps the sheet's default (unset) row height, same as rows 2-3.
$logs.Rows.Item(4).AutoFit()

# ---- Extend the conditional-formatting ranges from row 2:3 to 2:4 ---------
$colsToExtend = "D", "G", "H", "I", "J"
foreach ($col in $colsToExtend) {
    $oldRange = $logs.Range($col + "2:" + $col + "3")
    $newRange = $logs.Range($col + "2:" + $col + "4")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---- Dashboard sheet: recompute the category counts / ordering ------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A2").Value = "Overig"
$dash.Range("B2").Value = 2
$dash.Range("A3").Value = "Planning / Afspraak"
$dash.Range("B3").Value = 1
